$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E1").Value = "µA Keithley"

$keithley = @(29.1, 30, 35.9, 43.1, 66.6, 94, 98, 100.8, 106.8, 103, 106, 116, 133.5, 166, 271)
for ($i = 0; $i -lt $keithley.Length; $i++) {
    $ws.Cells.Item($i + 2, 5).Value = $keithley[$i]
}
